# Weekly data refresh: a new daily price record is inserted at the top of
# the data table (row 43), pushing the existing historical rows down by one.
#
# Before: data rows occupy 2..137 (row 1 is the header), dimension A1:R137.
# After : a new row is inserted at 43 (rows 43..137 shift to 44..138),
#         and the new row 43 receives the latest observation values,
#         giving dimension A1:R138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43; this shifts rows 43:137 down to
# 44:138 and carries their formatting/values with them automatically.
$row = $ws.Rows.Item(43)
$row.Insert()

# Populate the newly inserted row 43 with the new observation.
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value = 45274
$ws.Cells.Item(43, 5).Value = 15
$ws.Cells.Item(43, 6).Value = 100112040
$ws.Cells.Item(43, 7).Value = "Cilantro"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 450
$ws.Cells.Item(43, 11).Value = 1200
$ws.Cells.Item(43, 12).Value = 1500
$ws.Cells.Item(43, 13).Value = 1367
$ws.Cells.Item(43, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 684
$ws.Cells.Item(43, 17).Value = 2
$ws.Cells.Item(43, 18).Value = "Hortaliza"
